$d = $word.ActiveDocument

# The paragraph currently reads "Version 2." built from the runs:
#   "Versi" | "on" | " 2" | "."
# It needs to become "Version 1." built from the runs:
#   "Version" | " 1."
# (i.e. the trailing "." run is removed, the " 2" run becomes " 1.",
#  and the "Versi"/"on" runs are merged into a single "Version" run).

# Step 1: remove the trailing "." run entirely.
$d.Range(9, 10).Text = ""

# Step 2: turn the " 2" run into " 1."
$d.Range(7, 9).Text = " 1."

# Step 3: merge the "Versi" + "on" runs into a single "Version" run.
# A direct no-op rewrite ("Version" -> "Version") is ignored by the
# engine, so first force a real text change and then trim it back down,
# which causes the two runs to coalesce into one.
$d.Range(0, 7).Text = "Versionx"
$d.Range(7, 8).Text = ""
